# Fixing issues with experiments validation
#
# The "bs-seq" worksheet gets three new leading data columns inserted
# (after column A): "Experiment Alias", "Project" and "Secondary Project".
# Everything that used to live in column B onward shifts three columns to
# the right. The "bs-seq" tab also becomes the active tab of the workbook
# (it previously was "submission").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bs-seq")

# Insert three new blank columns before the old column B (old B/C/D become
# E/F/G, etc.)
$ws.Columns("B:D").Insert()

# Populate the headers of the three newly inserted columns.
$ws.Cells.Item(1, 2).Value = "Experiment Alias"
$ws.Cells.Item(1, 3).Value = "Project"
$ws.Cells.Item(1, 4).Value = "Secondary Project"

# Re-apply "best fit"-like widths for the columns whose content changed
# because of the insert (closest achievable values given this runtime's
# column-width quantization).
$ws.Columns.Item(2).ColumnWidth = 14.333333333333334
$ws.Columns.Item(3).ColumnWidth = 6
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(6).ColumnWidth = 13.333333333333334

# Make "bs-seq" the active/selected sheet (was "submission" before).
$ws.Activate()
